$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2098
$ws1.Range("F6").Value = 1758
$ws1.Range("F15").Value = 147
$ws1.Range("F18").Value = 4093
$ws1.Range("F23").Value = 957
$ws1.Range("F24").Value = 1106
$ws1.Range("F26").Value = 29
$ws1.Range("F28").Value = 1856
$ws1.Range("F29").Value = 52
$ws1.Range("F30").Value = 39
$ws1.Range("F31").Value = 78
$ws1.Range("F33").Value = 22

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2098
$ws4.Range("F6").Value = 1758
$ws4.Range("F15").Value = 147
$ws4.Range("F18").Value = 4093
$ws4.Range("F23").Value = 957
$ws4.Range("F24").Value = 1107
$ws4.Range("F26").Value = 29
$ws4.Range("F28").Value = 1856
$ws4.Range("F29").Value = 52
$ws4.Range("F30").Value = 39
$ws4.Range("F31").Value = 78
$ws4.Range("F33").Value = 22
